$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update description text for existing rows (row 2, 5, 6) ---
$ws.Range("K2").Value = "An inclusion criterion (rule) is_a *eligibility criterion* which defines and states a condition which, if met, makes an entity suitable for a given task or participation in a given process. For instance, in a study protocol, inclusion criteria indicate the conditions that prospective subjects MUST meet to be eligible for participation in a study."
$ws.Range("K5").Value = "That which constitutes a standard from which a judgment can be established, a reference point against which other things can be evaluated, or a basis for comparison."
$ws.Range("K6").Value = "Something distinguishable as an identifiable class based on common qualities"

# --- Add new semantics hyperlinks on existing rows (order matches rId allocation) ---
$ws.Hyperlinks.Add($ws.Range("J5"), "http://purl.obolibrary.org/obo/NCIT_C25466")
$ws.Range("J5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("J6"), "http://purl.obolibrary.org/obo/NCIT_C25284")
$ws.Range("J6").Style = "Hyperlink"

# --- Insert a new row 27 (organisations / type / string) ---
$ws.Rows.Item(27).Insert()

$ws.Range("A27").Value = "organisations"
$ws.Range("B27").Value = "ONTOLOGIES"
$ws.Range("D27").Value = "type"
$ws.Range("E27").Value = "string"
$ws.Range("K27").Value = "Something distinguishable as an identifiable class based on common qualities"

$ws.Hyperlinks.Add($ws.Range("J27"), "http://purl.obolibrary.org/obo/NCIT_C25284")
$ws.Range("J27").Style = "Hyperlink"

# --- Add new semantics hyperlink on row 2 (inclusionCriteria heading) ---
$ws.Hyperlinks.Add($ws.Range("J2"), "http://purl.obolibrary.org/obo/OBI_0500027")
$ws.Range("J2").Style = "Hyperlink"

# --- Add new semantics hyperlink + keep description on row 9 (dataproviders heading) ---
$ws.Hyperlinks.Add($ws.Range("J9"), "http://purl.obolibrary.org/obo/OBI_0000947")
$ws.Range("J9").Style = "Hyperlink"

# --- Row 28 (formerly row 27, pushed down by the insert) gains a new hyperlink cell ---
$ws.Hyperlinks.Add($ws.Range("J28"), "http://purl.obolibrary.org/obo/OBI_0000947")
$ws.Range("J28").Style = "Hyperlink"

# --- Update selection to match the final authored state ---
$ws.Range("J28").Select()
